$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 345
$ws.Range("J4").Value = 1500
$ws.Range("L4").Value = 1500
$ws.Range("N4").Value = -1728

$ws.Range("H11").Value = 3024.0833
$ws.Range("I11").Value = 3024.0833
$ws.Range("K11").Value = 3024.0833
$ws.Range("M11").Value = -2884.0833

$ws.Range("H19").Value = 4206.3335
$ws.Range("I19").Value = 3039.8572
$ws.Range("K19").Value = 3039.8572
$ws.Range("M19").Value = -2864.8572

$ws.Range("H101").Value = 4545753
$ws.Range("I101").Value = 11363983
$ws.Range("J101").Value = 266.5
$ws.Range("K101").Value = 34091949
$ws.Range("L101").Value = 799.5
$ws.Range("M101").Value = -34090327
$ws.Range("N101").Value = -4043.5

$ws.Range("H131").Value = 3715.4119
$ws.Range("I131").Value = 1810.8667
$ws.Range("K131").Value = 5432.6001
$ws.Range("M131").Value = -392.6000999999997

$ws.Range("H132").Value = 3423.1086
$ws.Range("I132").Value = 2913.7026
$ws.Range("J132").Value = 5517.3335
$ws.Range("K132").Value = 8741.1078
$ws.Range("L132").Value = 16552.0005
$ws.Range("M132").Value = -6211.1078
$ws.Range("N132").Value = -21612.0005

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws.Range("H135").Value = 993
$ws.Range("I135").Value = 992.72
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 8934.48
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -6399.48
$ws.Range("N135").Value = -14070

$ws.Range("H138").Value = 2522.5557
$ws.Range("I138").Value = 1042.875
$ws.Range("J138").Value = 3488.8774
$ws.Range("K138").Value = 3128.625
$ws.Range("L138").Value = 10466.6322
$ws.Range("M138").Value = 2011.375
$ws.Range("N138").Value = -20746.6322

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13932.825
$ws.Range("I32").Value = 8120.6665
$ws.Range("K32").Value = 8120.6665
$ws.Range("M32").Value = -7833.6665

$ws.Range("H122").Value = 4692.4
$ws.Range("I122").Value = 2410.3333
$ws.Range("J122").Value = 8115.5
$ws.Range("K122").Value = 7230.999899999999
$ws.Range("L122").Value = 24346.5
$ws.Range("M122").Value = -4780.999899999999
$ws.Range("N122").Value = -29246.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3684.92
$ws.Range("I86").Value = 3330.3333
$ws.Range("J86").Value = 5546.5
$ws.Range("K86").Value = 3330.3333
$ws.Range("L86").Value = 5546.5
$ws.Range("M86").Value = -2207.3333
$ws.Range("N86").Value = -7792.5

$ws.Range("H89").Value = 3684.92
$ws.Range("I89").Value = 3330.3333
$ws.Range("J89").Value = 5546.5
$ws.Range("K89").Value = 16651.6665
$ws.Range("L89").Value = 27732.5
$ws.Range("M89").Value = -11035.6665
$ws.Range("N89").Value = -38964.5

$ws.Range("H94").Value = 806.4211
$ws.Range("I94").Value = 672
$ws.Range("K94").Value = 672
$ws.Range("M94").Value = -221

$ws.Range("H99").Value = 4511.231
$ws.Range("I99").Value = 4107.1
$ws.Range("J99").Value = 4763.8125
$ws.Range("K99").Value = 4107.1
$ws.Range("L99").Value = 4763.8125
$ws.Range("M99").Value = -2609.1
$ws.Range("N99").Value = -7759.8125

$ws.Range("H134").Value = 3572.2173
$ws.Range("I134").Value = 2727.5293
$ws.Range("J134").Value = 5965.5
$ws.Range("K134").Value = 8182.5879
$ws.Range("L134").Value = 17896.5
$ws.Range("M134").Value = -5647.5879
$ws.Range("N134").Value = -22966.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 38367.145
$ws.Range("J68").Value = 38367.145
$ws.Range("L68").Value = 38367.145
$ws.Range("N68").Value = -39865.145

$ws.Range("H71").Value = 38367.145
$ws.Range("J71").Value = 38367.145
$ws.Range("L71").Value = 115101.435
$ws.Range("N71").Value = -122589.435

$ws.Range("H99").Value = 471192.38
$ws.Range("I99").Value = 912995.75
$ws.Range("K99").Value = 912995.75
$ws.Range("M99").Value = -911497.75

$ws.Range("H126").Value = 471192.38
$ws.Range("I126").Value = 912995.75
$ws.Range("K126").Value = 2738987.25
$ws.Range("M126").Value = -2736517.25

$ws.Range("H132").Value = 4218.9287
$ws.Range("I132").Value = 4883.2
$ws.Range("J132").Value = 3849.889
$ws.Range("K132").Value = 14649.6
$ws.Range("L132").Value = 11549.667
$ws.Range("M132").Value = -12119.6
$ws.Range("N132").Value = -16609.667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 574.3461
$ws.Range("I39").Value = 365
$ws.Range("J39").Value = 909.3
$ws.Range("K39").Value = 1095
$ws.Range("L39").Value = 2727.9
$ws.Range("M39").Value = -801
$ws.Range("N39").Value = -3315.9

$ws.Range("H113").Value = 649.5833
$ws.Range("J113").Value = 830.4286
$ws.Range("L113").Value = 2491.2858
$ws.Range("N113").Value = -6831.2858

$ws.Range("H131").Value = 10205650
$ws.Range("J131").Value = 1610.6666
$ws.Range("L131").Value = 4831.9998
$ws.Range("N131").Value = -14911.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 13516.056
$ws.Range("I122").Value = 11349.4
$ws.Range("K122").Value = 34048.2
$ws.Range("M122").Value = -31598.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1118.1714
$ws.Range("I22").Value = 900.2857
$ws.Range("J22").Value = 1263.4286
$ws.Range("K22").Value = 900.2857
$ws.Range("L22").Value = 1263.4286
$ws.Range("M22").Value = -605.2857
$ws.Range("N22").Value = -1853.4286

$ws.Range("H27").Value = 1118.1714
$ws.Range("I27").Value = 900.2857
$ws.Range("J27").Value = 1263.4286
$ws.Range("K27").Value = 900.2857
$ws.Range("L27").Value = 1263.4286
$ws.Range("M27").Value = -793.2857
$ws.Range("N27").Value = -1477.4286

$ws.Range("H122").Value = 2998.25
$ws.Range("I122").Value = 2998
$ws.Range("K122").Value = 8994
$ws.Range("M122").Value = -6544

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 9525.5
$ws.Range("I58").Value = 9525.5
$ws.Range("K58").Value = 9525.5
$ws.Range("M58").Value = -9217.5

$ws.Range("H59").Value = 11092
$ws.Range("I59").Value = 11092
$ws.Range("K59").Value = 11092
$ws.Range("M59").Value = -10354

$ws.Range("H62").Value = 9885.333000000001
$ws.Range("J62").Value = 10566.857
$ws.Range("L62").Value = 10566.857
$ws.Range("N62").Value = -11814.857

$ws.Range("H65").Value = 9885.333000000001
$ws.Range("J65").Value = 10566.857
$ws.Range("L65").Value = 52834.285
$ws.Range("N65").Value = -59074.285

$ws.Range("H98").Value = 22333.334
$ws.Range("J98").Value = 22333.334
$ws.Range("L98").Value = 22333.334
$ws.Range("N98").Value = -28323.334

$ws.Range("H122").Value = 9768.6
$ws.Range("I122").Value = 2889
$ws.Range("J122").Value = 33465
$ws.Range("K122").Value = 8667
$ws.Range("L122").Value = 100395
$ws.Range("M122").Value = -6217
$ws.Range("N122").Value = -105295

$ws.Range("H126").Value = 2846.4
$ws.Range("I126").Value = 2901
$ws.Range("J126").Value = 2696.25
$ws.Range("K126").Value = 8703
$ws.Range("L126").Value = 8088.75
$ws.Range("M126").Value = -6233
$ws.Range("N126").Value = -13028.75

$ws.Range("H132").Value = 2171.2856
$ws.Range("I132").Value = 2171.2856
$ws.Range("K132").Value = 6513.8568
$ws.Range("M132").Value = -3983.8568
